# Update the division-problem answers in the first (and only) table.
# The table has 5 "answer" rows (1, 5, 9, 13, 17) each holding 5 columns
# of "a÷b=c, d" text, separated by blank spacer rows. We target each
# cell directly by (row, column) so that duplicate source strings
# (e.g. "94÷7=13, 3" appearing twice) are updated independently.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text  = "63÷8=7, 7"
$t.Cell(1, 2).Range.Text  = "90÷9=10, 0"
$t.Cell(1, 3).Range.Text  = "97÷2=48, 1"
$t.Cell(1, 4).Range.Text  = "64÷3=21, 1"
$t.Cell(1, 5).Range.Text  = "23÷3=7, 2"

$t.Cell(5, 1).Range.Text  = "45÷5=9, 0"
$t.Cell(5, 2).Range.Text  = "21÷9=2, 3"
$t.Cell(5, 3).Range.Text  = "41÷2=20, 1"
$t.Cell(5, 4).Range.Text  = "32÷5=6, 2"
$t.Cell(5, 5).Range.Text  = "49÷9=5, 4"

$t.Cell(9, 1).Range.Text  = "83÷7=11, 6"
$t.Cell(9, 2).Range.Text  = "15÷2=7, 1"
$t.Cell(9, 3).Range.Text  = "63÷8=7, 7"
$t.Cell(9, 4).Range.Text  = "14÷2=7, 0"
$t.Cell(9, 5).Range.Text  = "66÷6=11, 0"

$t.Cell(13, 1).Range.Text = "39÷2=19, 1"
$t.Cell(13, 2).Range.Text = "47÷8=5, 7"
$t.Cell(13, 3).Range.Text = "69÷7=9, 6"
$t.Cell(13, 4).Range.Text = "48÷5=9, 3"
$t.Cell(13, 5).Range.Text = "66÷3=22, 0"

$t.Cell(17, 1).Range.Text = "12÷2=6, 0"
$t.Cell(17, 2).Range.Text = "69÷4=17, 1"
$t.Cell(17, 3).Range.Text = "63÷7=9, 0"
$t.Cell(17, 4).Range.Text = "53÷2=26, 1"
$t.Cell(17, 5).Range.Text = "43÷5=8, 3"
